# Add self feedback option and profile icon popup
# -> Appends two new "Pending" / "Q1" review rows (rows 6 and 7) that mirror
#    the existing pending-review rows, each with its own endDate/startDate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "yyyy-MM-dd HH:mm:ss"

# Row 6
$ws.Range("C6").Value = "Q1"
$ws.Range("D6").Value = 10
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 45311.229166666664
$ws.Range("G6").NumberFormat = $dateFormat
$ws.Range("J6").Value = 6
$ws.Range("L6").Value = "Pending"
$ws.Range("M6").Value = 45308.229166666664
$ws.Range("M6").NumberFormat = $dateFormat

# Row 7
$ws.Range("C7").Value = "Q1"
$ws.Range("D7").Value = 10
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 45311.229166666664
$ws.Range("G7").NumberFormat = $dateFormat
$ws.Range("J7").Value = 6
$ws.Range("L7").Value = "Pending"
$ws.Range("M7").Value = 45308.229166666664
$ws.Range("M7").NumberFormat = $dateFormat
